$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the week 4 hours for the first team member
$ws.Range("B7").Value = 10

# Recalculate so dependent formulas (C7:C13 running totals) update
$excel.Calculate()

# Reflect the final cell selection made by the user after entering the value
$ws.Range("B8").Select()
